$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.150.41"
$ws.Range("E2").Value = "  -4.50%  "
$ws.Range("D3").Value = "'1.655.95"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("D5").Value = "'217.86"
$ws.Range("E5").Value = "  -2.91%  "
$ws.Range("D6").Value = "'0.5177"
$ws.Range("E6").Value = "  -2.74%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.06451"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").Value = "'0.2572"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("D10").Value = "'19.89"
$ws.Range("E10").Value = "  -5.11%  "
$ws.Range("D11").Value = "'0.07787"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "'1.669.52"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.294"
$ws.Range("E13").Value = "  -5.67%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.882.92"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "'0.5533"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "'0.0₅8052"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'64.33"
$ws.Range("E17").Value = "  -5.13%  "
$ws.Range("D18").Value = "'26.184.68"
$ws.Range("E18").Value = "  -4.32%  "
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'211.28"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").Value = "'4.388"
$ws.Range("E21").Value = "  -5.82%  "
$ws.Range("D22").Value = "'10.07"
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "'5.913"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'143.90"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'1.763"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("D28").Value = "'6.974"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "'15.75"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").Value = "'0.05281"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").Value = "'1.252"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").Value = "'3.363"
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("D33").Value = "'3.233"
$ws.Range("E33").Value = "  -5.65%  "
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("D35").Value = "'2.765"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "'2.369"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'0.9238"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").Value = "'1.165.25"
$ws.Range("E38").Value = "  +11.26%  "
$ws.Range("D39").Value = "'0.5689"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").Value = "'0.01591"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "'0.8371"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'5.649"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").Value = "'99.87"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "'1.793.49"
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("E46").Value = "  -6.80%  "
$ws.Range("D47").Value = "'0.4510"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'55.94"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'7.868"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").Value = "'0.05074"
$ws.Range("E51").Value = "  -2.97%  "

Write-Output "edit applied"
